$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column data-type headers
$ws.Range("D1").Value = "varchar(50)"
$ws.Range("E1").Value = "varchar(50)"
$ws.Range("P1").Value = "int(255)"
$ws.Range("Q1").Value = "int(255)"
$ws.Range("R1").Value = "varchar(6)"

# Row 2: two new columns (D, E) are introduced before the old "기초재고수량"
# column, shifting the remaining labels from D..Q over to F..R, plus a
# couple of relabeled cells (A2, C2) and a new note on R2. Write the final
# label for every column directly.
$ws.Range("A2").Value = "ID (입력x)"
$ws.Range("C2").Value = "버젼코드"
$ws.Range("D2").Value = "계정코드"
$ws.Range("E2").Value = "완제품코드"
$ws.Range("F2").Value = "기초재고수량"
$ws.Range("G2").Value = "기초재고금액"
$ws.Range("H2").Value = "생산입고수량"
$ws.Range("I2").Value = "생산입고금액"
$ws.Range("J2").Value = "판매출고수량"
$ws.Range("K2").Value = "판매출고금액"
$ws.Range("L2").Value = "LOSS출고수량"
$ws.Range("M2").Value = "LOSS출고금액"
$ws.Range("N2").Value = "개발출고수량"
$ws.Range("O2").Value = "개발출고금액"
$ws.Range("P2").Value = "기말재고수량"
$ws.Range("Q2").Value = "기말재고금액"
$ws.Range("R2").Value = "년월 ex) 200001"
